$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 131.16667
$ws.Range("I2").Value = 59.25
$ws.Range("J2").Value = 275
$ws.Range("K2").Value = 59.25
$ws.Range("L2").Value = 275
$ws.Range("M2").Value = 53.75
$ws.Range("N2").Value = -501
$ws.Range("H4").Value = 463.64285
$ws.Range("I4").Value = 408.45456
$ws.Range("J4").Value = 666
$ws.Range("K4").Value = 408.45456
$ws.Range("L4").Value = 666
$ws.Range("M4").Value = -294.45456
$ws.Range("N4").Value = -894
$ws.Range("H15").Value = 2362.4602
$ws.Range("I15").Value = 2362.4602
$ws.Range("K15").Value = 7087.3806
$ws.Range("M15").Value = -6918.3806
$ws.Range("H33").Value = 879.2222
$ws.Range("I33").Value = 909
$ws.Range("J33").Value = 641
$ws.Range("K33").Value = 909
$ws.Range("L33").Value = 641
$ws.Range("M33").Value = -680
$ws.Range("N33").Value = -1099
$ws.Range("H86").Value = 45074.07
$ws.Range("I86").Value = 62009.7
$ws.Range("J86").Value = 2735
$ws.Range("K86").Value = 62009.7
$ws.Range("L86").Value = 2735
$ws.Range("M86").Value = -60886.7
$ws.Range("N86").Value = -4981
$ws.Range("H89").Value = 45074.07
$ws.Range("I89").Value = 62009.7
$ws.Range("J89").Value = 2735
$ws.Range("K89").Value = 310048.5
$ws.Range("L89").Value = 13675
$ws.Range("M89").Value = -304432.5
$ws.Range("N89").Value = -24907
$ws.Range("H132").Value = 9733.759
$ws.Range("I132").Value = 7115.2
$ws.Range("J132").Value = 26099.75
$ws.Range("K132").Value = 21345.6
$ws.Range("L132").Value = 78299.25
$ws.Range("M132").Value = -18815.6
$ws.Range("N132").Value = -83359.25
$ws.Range("H137").Value = 16131383
$ws.Range("I137").Value = 19233192
$ws.Range("J137").Value = 1978
$ws.Range("K137").Value = 57699576
$ws.Range("L137").Value = 5934
$ws.Range("M137").Value = -57697026
$ws.Range("N137").Value = -11034
$ws.Range("H138").Value = 1794.4557
$ws.Range("I138").Value = 1347.6
$ws.Range("J138").Value = 2385.8823
$ws.Range("K138").Value = 4042.8
$ws.Range("L138").Value = 7157.646900000001
$ws.Range("M138").Value = 1097.2
$ws.Range("N138").Value = -17437.6469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2546.35
$ws.Range("I122").Value = 1301.5555
$ws.Range("K122").Value = 3904.6665
$ws.Range("M122").Value = -1454.6665
$ws.Range("H132").Value = 2086864.5
$ws.Range("I132").Value = 1510.5
$ws.Range("J132").Value = 5214895.5
$ws.Range("K132").Value = 4531.5
$ws.Range("L132").Value = 15644686.5
$ws.Range("M132").Value = -2001.5
$ws.Range("N132").Value = -15649746.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 666.75
$ws.Range("I94").Value = 680.6667
$ws.Range("J94").Value = 625
$ws.Range("K94").Value = 680.6667
$ws.Range("L94").Value = 625
$ws.Range("M94").Value = -229.6667
$ws.Range("N94").Value = -1527

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 7727.769
$ws.Range("I7").Value = 10023.8
$ws.Range("J7").Value = 74.333336
$ws.Range("K7").Value = 10023.8
$ws.Range("L7").Value = 74.333336
$ws.Range("M7").Value = -9910.799999999999
$ws.Range("N7").Value = -300.333336
$ws.Range("H22").Value = 340
$ws.Range("I22").Value = 245
$ws.Range("J22").Value = 593.3333
$ws.Range("K22").Value = 245
$ws.Range("L22").Value = 593.3333
$ws.Range("M22").Value = 105
$ws.Range("N22").Value = -1293.3333
$ws.Range("H86").Value = 38465884
$ws.Range("I86").Value = 71431310
$ws.Range("J86").Value = 6217.3335
$ws.Range("K86").Value = 71431310
$ws.Range("L86").Value = 6217.3335
$ws.Range("M86").Value = -71430187
$ws.Range("N86").Value = -8463.333500000001
$ws.Range("H89").Value = 38465884
$ws.Range("I89").Value = 71431310
$ws.Range("J89").Value = 6217.3335
$ws.Range("K89").Value = 357156550
$ws.Range("L89").Value = 31086.6675
$ws.Range("M89").Value = -357150934
$ws.Range("N89").Value = -42318.6675
$ws.Range("H105").Value = 417717.75
$ws.Range("I105").Value = 589272.9
$ws.Range("K105").Value = 589272.9
$ws.Range("M105").Value = -587525.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 71.416664
$ws.Range("I2").Value = 40.42857
$ws.Range("J2").Value = 114.8
$ws.Range("K2").Value = 242.57142
$ws.Range("L2").Value = 688.8
$ws.Range("M2").Value = -129.57142
$ws.Range("N2").Value = -914.8
$ws.Range("H14").Value = 168.14285
$ws.Range("I14").Value = 168.14285
$ws.Range("K14").Value = 504.42855
$ws.Range("M14").Value = -331.42855
$ws.Range("H68").Value = 856.73
$ws.Range("I68").Value = 721.68335
$ws.Range("J68").Value = 1059.3
$ws.Range("K68").Value = 2165.05005
$ws.Range("L68").Value = 3177.9
$ws.Range("M68").Value = -1354.05005
$ws.Range("N68").Value = -4799.9
$ws.Range("H71").Value = 856.73
$ws.Range("I71").Value = 721.68335
$ws.Range("J71").Value = 1059.3
$ws.Range("K71").Value = 6495.15015
$ws.Range("L71").Value = 9533.699999999999
$ws.Range("M71").Value = -2439.15015
$ws.Range("N71").Value = -17645.7
$ws.Range("H105").Value = 6109.231
$ws.Range("I105").Value = 3000
$ws.Range("J105").Value = 6368.3335
$ws.Range("K105").Value = 9000
$ws.Range("L105").Value = 19105.0005
$ws.Range("M105").Value = -6379
$ws.Range("N105").Value = -24347.0005
$ws.Range("H107").Value = 9455236
$ws.Range("I107").Value = 31251652
$ws.Range("J107").Value = 277798.78
$ws.Range("K107").Value = 93754956
$ws.Range("L107").Value = 833396.3400000001
$ws.Range("M107").Value = -93753036
$ws.Range("N107").Value = -837236.3400000001
$ws.Range("H131").Value = 481627.53
$ws.Range("I131").Value = 532.3077
$ws.Range("J131").Value = 962722.75
$ws.Range("K131").Value = 1596.9231
$ws.Range("L131").Value = 2888168.25
$ws.Range("M131").Value = 3443.0769
$ws.Range("N131").Value = -2898248.25
$ws.Range("H139").Value = 2775.7585
$ws.Range("I139").Value = 957.2
$ws.Range("J139").Value = 6817
$ws.Range("K139").Value = 2871.6
$ws.Range("L139").Value = 20451
$ws.Range("M139").Value = 2268.4
$ws.Range("N139").Value = -30731

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5431.6665
$ws.Range("I70").Value = 5940
$ws.Range("K70").Value = 5940
$ws.Range("M70").Value = -5670
$ws.Range("H73").Value = 5431.6665
$ws.Range("I73").Value = 5940
$ws.Range("K73").Value = 5940
$ws.Range("M73").Value = -5004
$ws.Range("H122").Value = 3776.1538
$ws.Range("J122").Value = 5460
$ws.Range("L122").Value = 16380
$ws.Range("N122").Value = -21280
$ws.Range("H132").Value = 888260.8
$ws.Range("I132").Value = 1390353.6
$ws.Range("K132").Value = 4171060.8
$ws.Range("M132").Value = -4168530.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 19929.875
$ws.Range("I122").Value = 23995.8
$ws.Range("K122").Value = 71987.39999999999
$ws.Range("M122").Value = -69537.39999999999
